{"js": "const body = context.document.body;\n\n// 1) Split the paragraph \"TODO: constraints test doesn't work\" into two\n//    paragraphs: \"TODO: \" and \"constraints test doesn't work\".\nconst todoParaResults = body.search(\"TODO: constraints test doesn\\u2019t work\", { matchCase: true });\nawait context.sync();\nconst todoPara = todoParaResults.items[0];\n\nconst todoPrefixResults = todoPara.search(\"TODO: \", { matchCase: true });\nawait context.sync();\nconst todoPrefix = todoPrefixResults.items[0];\ntodoPrefix.insertText(\"\\r\", \"End\");\nawait context.sync();\n\n// 2) Change the text \"Select mode is broken\" (keeping its paragraph /\n//    bookmark intact) to \"Make all things work as java executables\".\nconst selectModeResults = body.search(\"Select mode is broken\", { matchCase: true });\nawait context.sync();\nconst selectModeRange = selectModeResults.items[0];\nselectModeRange.insertText(\"Make all things work as java executables\", \"Replace\");\nawait context.sync();\n\n// 3) Add two new paragraphs right after that one.\nconst paragraphs1 = body.paragraphs;\nparagraphs1.load(\"items\");\nawait context.sync();\nconst lastPara1 = paragraphs1.items[paragraphs1.items.length - 1];\nlastPara1.insertParagraph(\"Test node-edge editor\", \"After\");\nawait context.sync();\n\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nconst lastPara2 = paragraphs2.items[paragraphs2.items.length - 1];\nlastPara2.insertParagraph(\"Make a cool demo that uses processing somehow\", \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rsquo = [char]0x2019\n\n# 1) Split the paragraph \"TODO: constraints test doesn't work\" into two\n#    paragraphs: \"TODO: \" and \"constraints test doesn't work\".\n$anchor = $d.Content.Duplicate\n$anchor.Find.Text = \"TODO: constraints test doesn\" + $rsquo + \"t work\"\n$anchor.Find.Execute() | Out-Null\n\n$splitPoint = $anchor.Duplicate\n$splitPoint.Find.Text = \"TODO: \"\n$splitPoint.Find.Execute() | Out-Null\n$splitPoint.InsertParagraphAfter()\n\n# 2) Change the text \"Select mode is broken\" (keeping its paragraph /\n#    bookmark intact) to \"Make all things work as java executables\".\n$selectMode = $d.Content.Duplicate\n$selectMode.Find.Text = \"Select mode is broken\"\n$selectMode.Find.Execute() | Out-Null\n$selectMode.Text = \"Make all things work as java executables\"\n\n# 3) Add two new paragraphs right after that one.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Test node-edge editor\"\n\n$lastPara2 = $d.Paragraphs.Last\n$lastPara2.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Make a cool demo that uses processing somehow\"\n"}
